# Apply the crypto price/volume refresh for Mon Feb 12 21:24:44 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume(1h)) hold formatted strings (thousand-dot-separated
# prices, space-padded percentages) that must stay text, not be re-parsed as
# numbers/dates by Excel. Force Text format for the edited range first, then
# restore the original (General) number format afterwards by pasting formats
# from the untouched B column (same rows) so no stray styling is introduced.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "50.180.60"
$ws.Range("E2").Value = "  +4.44%  "
$ws.Range("D3").Value = "2.645.98"
$ws.Range("E3").Value = "  +5.96%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "328.42"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").Value = "111.34"
$ws.Range("E6").Value = "  +3.57%  "
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.562"
$ws.Range("E9").Value = "  +4.56%  "
$ws.Range("E10").Value = "  +3.54%  "
$ws.Range("D11").Value = "20.70"
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("D15").Value = "3.059.59"
$ws.Range("E15").Value = "  +5.95%  "
$ws.Range("D16").Value = "2.636.39"
$ws.Range("E16").Value = "  +5.57%  "
$ws.Range("D17").Value = "0.883"
$ws.Range("E17").Value = "  +5.59%  "
$ws.Range("D18").Value = "50.104.95"
$ws.Range("E18").Value = "  +4.55%  "
$ws.Range("E19").Value = "  +11.82%  "
$ws.Range("D20").Value = "13.36"
$ws.Range("E20").Value = "  +3.37%  "
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  +2.68%  "
$ws.Range("D23").Value = "73.09"
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("D24").Value = "279.99"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("E25").Value = "  +2.88%  "
$ws.Range("D26").Value = "26.72"
$ws.Range("E26").Value = "  +4.36%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.23"
$ws.Range("E28").Value = "  +6.40%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "36.70"
$ws.Range("E29").Value = "  +5.01%  "
$ws.Range("D30").Value = "9.96"
$ws.Range("E30").Value = "  +2.23%  "
$ws.Range("D31").Value = "0.145"
$ws.Range("E31").Value = "  +3.38%  "
$ws.Range("D32").Value = "49.91"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").Value = "19.84"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("E34").Value = "  +3.06%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "0.0796"
$ws.Range("E36").Value = "  +2.14%  "
$ws.Range("E37").Value = "  +6.96%  "
$ws.Range("E38").Value = "  +3.31%  "
$ws.Range("E39").Value = "  +7.75%  "
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").Value = "123.38"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").Value = "22.69"
$ws.Range("E42").Value = "  +7.22%  "
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").Value = "0.0315"
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("D45").Value = "3.35"
$ws.Range("E45").Value = "  +6.58%  "
$ws.Range("D46").Value = "2.065.39"
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("D47").Value = "2.31"
$ws.Range("E47").Value = "  +15.91%  "
$ws.Range("E48").Value = "  +8.96%  "
$ws.Range("D49").Value = "9.06"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").Value = "5.40"
$ws.Range("E50").Value = "  +4.61%  "
$ws.Range("D51").Value = "81.90"
$ws.Range("E51").Value = "  +2.18%  "

# Restore D:E to their original (General) formatting/style.
$ws.Range("B2:B51").Copy()
$ws.Range("D2:D51").PasteSpecial(-4122)
$ws.Range("B2:B51").Copy()
$ws.Range("E2:E51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
